{"js": "const replacements = [\n  [\"71\u00d726=1846\", \"50\u00d739=1950\"],\n  [\"88\u00d711=968\", \"62\u00d767=4154\"],\n  [\"38\u00d794=3572\", \"74\u00d756=4144\"],\n  [\"35\u00d763=2205\", \"90\u00d782=7380\"],\n  [\"85\u00d747=3995\", \"33\u00d727=891\"],\n  [\"91\u00d717=1547\", \"93\u00d719=1767\"],\n  [\"27\u00d745=1215\", \"53\u00d782=4346\"],\n  [\"94\u00d772=6768\", \"81\u00d760=4860\"],\n  [\"13\u00d787=1131\", \"52\u00d732=1664\"],\n  [\"62\u00d773=4526\", \"41\u00d726=1066\"],\n  [\"85\u00d763=5355\", \"30\u00d786=2580\"],\n  [\"50\u00d794=4700\", \"58\u00d787=5046\"],\n  [\"80\u00d771=5680\", \"31\u00d727=837\"],\n  [\"50\u00d736=1800\", \"17\u00d775=1275\"],\n  [\"37\u00d776=2812\", \"73\u00d746=3358\"],\n  [\"16\u00d767=1072\", \"47\u00d755=2585\"],\n  [\"20\u00d714=280\", \"51\u00d754=2754\"],\n  [\"43\u00d722=946\", \"81\u00d723=1863\"],\n  [\"69\u00d783=5727\", \"95\u00d731=2945\"],\n  [\"72\u00d752=3744\", \"29\u00d798=2842\"],\n  [\"63\u00d711=693\", \"77\u00d737=2849\"],\n  [\"97\u00d782=7954\", \"36\u00d796=3456\"],\n  [\"67\u00d720=1340\", \"41\u00d721=861\"],\n  [\"95\u00d777=7315\", \"70\u00d732=2240\"],\n  [\"14\u00d784=1176\", \"54\u00d792=4968\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"71\u00d726=1846\"; new=\"50\u00d739=1950\"},\n  @{old=\"88\u00d711=968\"; new=\"62\u00d767=4154\"},\n  @{old=\"38\u00d794=3572\"; new=\"74\u00d756=4144\"},\n  @{old=\"35\u00d763=2205\"; new=\"90\u00d782=7380\"},\n  @{old=\"85\u00d747=3995\"; new=\"33\u00d727=891\"},\n  @{old=\"91\u00d717=1547\"; new=\"93\u00d719=1767\"},\n  @{old=\"27\u00d745=1215\"; new=\"53\u00d782=4346\"},\n  @{old=\"94\u00d772=6768\"; new=\"81\u00d760=4860\"},\n  @{old=\"13\u00d787=1131\"; new=\"52\u00d732=1664\"},\n  @{old=\"62\u00d773=4526\"; new=\"41\u00d726=1066\"},\n  @{old=\"85\u00d763=5355\"; new=\"30\u00d786=2580\"},\n  @{old=\"50\u00d794=4700\"; new=\"58\u00d787=5046\"},\n  @{old=\"80\u00d771=5680\"; new=\"31\u00d727=837\"},\n  @{old=\"50\u00d736=1800\"; new=\"17\u00d775=1275\"},\n  @{old=\"37\u00d776=2812\"; new=\"73\u00d746=3358\"},\n  @{old=\"16\u00d767=1072\"; new=\"47\u00d755=2585\"},\n  @{old=\"20\u00d714=280\"; new=\"51\u00d754=2754\"},\n  @{old=\"43\u00d722=946\"; new=\"81\u00d723=1863\"},\n  @{old=\"69\u00d783=5727\"; new=\"95\u00d731=2945\"},\n  @{old=\"72\u00d752=3744\"; new=\"29\u00d798=2842\"},\n  @{old=\"63\u00d711=693\"; new=\"77\u00d737=2849\"},\n  @{old=\"97\u00d782=7954\"; new=\"36\u00d796=3456\"},\n  @{old=\"67\u00d720=1340\"; new=\"41\u00d721=861\"},\n  @{old=\"95\u00d777=7315\"; new=\"70\u00d732=2240\"},\n  @{old=\"14\u00d784=1176\"; new=\"54\u00d792=4968\"},\n)\n\nforeach ($pair in $pairs) {\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n}\n"}
